$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replaceText, 2) | Out-Null
}

# 1. "significantly impacts" -> "significantly impact"
Replace-Text "significantly impacts" "significantly impact"

# 2. "longitudinal assessments" -> "longitudinal assessment"
Replace-Text "longitudinal assessments" "longitudinal assessment"

# 3. "(18+)" -> "(18+ years)"
Replace-Text "(18+)" "(18+ years)"

# 4. "two factors unique to the" -> "two novel factors unique to the"
Replace-Text "two factors unique to the" "two novel factors unique to the"

# 5. "environment have with adult" -> "environment may have with adult"
Replace-Text "environment have with adult" "environment may have with adult"

# 6. test-retest reliability sentence
Replace-Text "has an excellent overall test-retest reliability." "has excellent overall test-retest reliability (General Factor ICC = 0.84)."

# 7. Replace full "Along with validating..." sentence
Replace-Text "Along with validating this scale, we also discuss results that demonstrate the utility of measuring the childhood home musical environment in adults. Highlights of these findings are as follows:" "We also present demonstrating the utility of measuring the childhood home musical environment in adults. Highlights of these findings are as follows:"

# 8. First bullet restructure
Replace-Text "Scores on the Caregiver Beliefs and Child Engagement with Music factors, as well as overall " "Caregiver Beliefs, Child Engagement with Music, and overall "

# 9. melodic perception sentence
Replace-Text "were positively associated to performance on a melodic perception task" "were significantly positively associated with performance on a melodic perception task"

# 10. second bullet restructure: Overall -> Caregiver Singing, Attitude..., and overall
Replace-Text "Overall Music@Home – Retrospective scores, as well as Caregiver Singing and Attitude towards Childhood Home Musical Environment, were negatively correlated with adult trait level anxiety" "Caregiver Singing, Attitude towards Childhood Home Musical Environment, and overall Music@Home – Retrospective scores were significantly negatively correlated with adult trait level anxiety"

# 11. third bullet restructure: resilience sentence
Replace-Text "All scores (general and factor) on the Music@Home – Retrospective scale were positively associated with adult resilience" "All Music@Home – Retrospective scores (i.e., factor and general scores) were significantly positively associated with adult self-reported resilience"

# 12. Double space fix: ",  Goldsmiths" -> ", Goldsmiths"
Replace-Text ",  Goldsmiths" ", Goldsmiths"
